$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1) Get rid of the rows that won't exist anymore (old rows 4-6),
#    and strip the bold/border/centered header-style that used to
#    live on column A of the rows that survive (old rows 2 & 3),
#    since in the new layout those become plain (unstyled) data
#    cells. Do this BEFORE copying any formatting around, because
#    any Clear-type call wipes the paste buffer.
# ------------------------------------------------------------------
$ws.Range("A4:D6").Clear()
$ws.Range("A2:A3").ClearFormats()

# ------------------------------------------------------------------
# 2) Grab the existing header style (bold, thin border, centered /
#    top-aligned) that already lives on B1 so we can stamp it onto
#    the full new header row without registering brand-new style
#    entries.
# ------------------------------------------------------------------
$ws.Range("B1").Copy()

# ------------------------------------------------------------------
# 3) Write the new header row (A1:R1).
# ------------------------------------------------------------------
$headers = @( `
  "stocks", `
  "predictions_1", `
  "predictions_2", `
  "predictions_3", `
  "predictions_4", `
  "predictions_5", `
  "predictions_6", `
  "predictions_7", `
  "growth_index", `
  "analisis r2", `
  "ultimo_preco", `
  "delta_ultimo_preco_vs_1_prediction", `
  "delta_1_prediction_vs_2_prediction", `
  "delta_2_prediction_vs_3_prediction", `
  "delta_3_prediction_vs_4_prediction", `
  "delta_4_prediction_vs_5_prediction", `
  "delta_5_prediction_vs_6_prediction", `
  "delta_6_prediction_vs_7_prediction" `
)

for ($i = 0; $i -lt $headers.Length; $i++) {
  $ws.Cells.Item(1, $i + 1).Value = $headers[$i]
}

# ------------------------------------------------------------------
# 4) Write the two data rows (ITSA4.SA / PETR3.SA).
# ------------------------------------------------------------------
$row2 = @(10.45422884127583, 11.22809856623553, 12.0867511430697, 12.75533252936839, 13.19144673365736, 13.44083089038486, 13.56252942915853, 29.73247128100325, 79.33427948756794, 10.43000030517578, 0.2322966001067739, 7.402456333309582, 7.647355175668458, 5.531522725873628, 3.419073577931764, 1.89049891010975, 0.9054391039227339)
$row3 = @(37.87024312050949, 41.34436253349327, 43.26341375248172, 44.28320988827728, 44.85279995381189, 45.14278428635863, 45.26521323270936, 19.52712605690923, 82.52300908964297, 36.97999954223633, 2.407365033242859, 9.173744678449891, 4.641627301506501, 2.357179074286675, 1.286243853983571, 0.6465244819618077, 0.2712037998677808)

$ws.Cells.Item(2, 1).Value = "ITSA4.SA"
for ($i = 0; $i -lt $row2.Length; $i++) {
  $ws.Cells.Item(2, $i + 2).Value = $row2[$i]
}

$ws.Cells.Item(3, 1).Value = "PETR3.SA"
for ($i = 0; $i -lt $row3.Length; $i++) {
  $ws.Cells.Item(3, $i + 2).Value = $row3[$i]
}

# ------------------------------------------------------------------
# 5) Stamp the header style across the whole new header row.
# ------------------------------------------------------------------
$ws.Range("A1:R1").PasteSpecial(-4122)
